# Fruta / hortaliza, semanal
# A new weekly record is added at row 28 (date 44459 = 2021-09-20),
# pushing all the existing records (rows 28-125) down by one row,
# so the former last record (row 125) becomes row 126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28; this shifts rows 28:125 down to 29:126
$ws.Rows(28).Insert()

# Populate the new row 28 with the new weekly record
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 44459
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100104
$ws.Range("H28").Value = "Frutos de pepita"
$ws.Range("I28").Value = 100104005
$ws.Range("J28").Value = "Pera"
$ws.Range("K28").Value = "Packham's Triumph"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 9500
$ws.Range("Q28").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R28").Value = "Provincia de Curicó"
$ws.Range("S28").Value = 594
$ws.Range("T28").Value = 16

# Make sure the date cell keeps the same date number format as the rest of column D
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
